$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$xlPasteFormats = -4122

# --- Data edits: refresh the QR-code test fixtures ---
# (numeroCuenta / usuario / clave / numeroDocumento columns, rows 2-5)

# Stash a copy of the original formatting for columns B, D and E (rows 2-5 all
# share the same per-column style) in a scratch area so it can be restored
# after the .Value assignment below, which otherwise resets the cell style.
$ws.Range("B2").Copy()
$ws.Range("AA1").PasteSpecial($xlPasteFormats)
$ws.Range("D2").Copy()
$ws.Range("AA2").PasteSpecial($xlPasteFormats)
$ws.Range("E2").Copy()
$ws.Range("AA3").PasteSpecial($xlPasteFormats)

# numeroCuenta (column Q) - each row gets its own new account number
$ws.Range("Q2").Value = "406-139440-02"
$ws.Range("Q3").Value = "406-739440-04"
$ws.Range("Q4").Value = "406-739440-03"
$ws.Range("Q5").Value = "406-739440-03"

# usuario (column D)
$ws.Range("D2:D5").Value = "userrobot9"

# clave (column E)
$ws.Range("E2:E5").Value = "6789"

# numeroDocumento (column B)
$ws.Range("B2:B5").Value = "22493944"

foreach ($row in 2..5) {
    $ws.Range("AA1").Copy()
    $ws.Range("B$row").PasteSpecial($xlPasteFormats)

    $ws.Range("AA2").Copy()
    $ws.Range("D$row").PasteSpecial($xlPasteFormats)

    $ws.Range("AA3").Copy()
    $ws.Range("E$row").PasteSpecial($xlPasteFormats)
}

# Clean up the scratch cells used to stash formatting
$ws.Range("AA1:AA3").Clear()

# --- View state: scroll back to the top-left and move the active selection to F9 ---
$win = $excel.ActiveWindow
$win.ScrollColumn = 1
$win.ScrollRow = 1
$ws.Range("F9").Select()
